$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.951.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.01%  "

$ws.Range("D3").Value = "'2.317.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.33%  "

$ws.Range("D5").Value = "'96.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.54%  "

$ws.Range("D6").Value = "'272.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.30%  "

$ws.Range("D7").Value = "'0.629"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  -2.19%  "

$ws.Range("D10").Value = "'45.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.05%  "

$ws.Range("D11").Value = "'0.0954"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.67%  "

$ws.Range("D12").Value = "'8.00"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.83%  "

$ws.Range("D13").Value = "'0.106"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.38%  "

$ws.Range("D14").Value = "'2.655.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.99%  "

$ws.Range("D15").Value = "'15.45"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.69%  "

$ws.Range("E16").Value = "  +6.65%  "

$ws.Range("D17").Value = "'2.319.58"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.19%  "

$ws.Range("D18").Value = "'43.867.32"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.08%  "

$ws.Range("E19").Value = "  +4.04%  "

$ws.Range("D20").Value = "'6.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.91%  "

$ws.Range("D21").Value = "'73.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.86%  "

$ws.Range("D22").Value = "'240.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.04%  "

$ws.Range("E23").Value = "  -2.03%  "

$ws.Range("D24").Value = "'9.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.84%  "

$ws.Range("E26").Value = "  +0.75%  "

$ws.Range("E27").Value = "  -0.95%  "

$ws.Range("D28").Value = "'3.50"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.16%  "

$ws.Range("D29").Value = "'2.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.97%  "

$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "'22.45"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.69%  "

$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").Value = "'38.07"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.75%  "

$ws.Range("D32").Value = "'175.56"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.66%  "

$ws.Range("D33").Value = "'0.0912"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.05%  "

$ws.Range("D34").Value = "'5.49"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.55%  "

$ws.Range("E35").Value = "  +1.94%  "

$ws.Range("D36").Value = "'0.0365"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.38%  "

$ws.Range("E37").Value = "  -3.40%  "

$ws.Range("D38").Value = "'4.45"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.65%  "

$ws.Range("D39").Value = "'3.36"
$ws.Range("D39").Style = "Normal"

$ws.Range("E40").Value = "  +7.35%  "

$ws.Range("D41").Value = "'2.38"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +9.81%  "

$ws.Range("E42").Value = "  +22.41%  "

$ws.Range("D43").Value = "'12.40"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.36%  "

$ws.Range("D44").Value = "'62.92"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.16%  "

$ws.Range("D45").Value = "'9.14"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.94%  "

$ws.Range("D46").Value = "'5.32"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.95%  "

$ws.Range("E47").Value = "  +3.80%  "

$ws.Range("D48").Value = "'100.44"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.87%  "

$ws.Range("E49").Value = "  +0.10%  "

$ws.Range("D50").Value = "'0.193"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +16.92%  "

$ws.Range("D51").Value = "'2.542.27"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.25%  "
